$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 17 with data, matching the pattern of existing rows
# Copy A16's formatting (date style) down to A17 first, then set its value
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Cells.Item(17, 1).Value = 42622.887291666666

$ws.Cells.Item(17, 2).Value = -22
$ws.Cells.Item(17, 3).Value = 63
$ws.Cells.Item(17, 4).Value = 33
$ws.Cells.Item(17, 5).Value = 63
$ws.Cells.Item(17, 6).Value = 18
$ws.Cells.Item(17, 7).Value = 8723
$ws.Cells.Item(17, 8).Value = 15227
$ws.Cells.Item(17, 9).Value = 1632
$ws.Cells.Item(17, 10).Value = 206
$ws.Cells.Item(17, 11).Value = 107
$ws.Cells.Item(17, 12).Value = 18
$ws.Cells.Item(17, 13).Value = 4
$ws.Cells.Item(17, 14).Value = "Bag"

$wb.Save()
